# Chromatin model workbook update:
#  - renumber the "#" (reaction index) column on the Reactions sheet to
#    account for the newly-added methylated-histone reactions
#  - set the methylation rate constants (E113:E127) to 0.001 with a
#    three-decimal number format
#  - restore the selection state on the Species and Reactions sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Species sheet: just a selection change
# ---------------------------------------------------------------------
$wsSpecies = $wb.Worksheets.Item("Species")
$wsSpecies.Activate()
$wsSpecies.Range("C16").Select()

# ---------------------------------------------------------------------
# Reactions sheet: renumber column A and update the methylation rates
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Reactions")
$ws.Activate()

$ws.Range("A58").Value = 29
$ws.Range("A61").Value = 30
$ws.Range("A62").Value = 31
$ws.Range("A65").Value = 32
$ws.Range("A66").Value = 33
$ws.Range("A69").Value = 34
$ws.Range("A70").Value = 35
$ws.Range("A73").Value = 36
$ws.Range("A74").Value = 37
$ws.Range("A77").Value = 38
$ws.Range("A78").Value = 39
$ws.Range("A81").Value = 40
$ws.Range("A82").Value = 41
$ws.Range("A85").Value = 42
$ws.Range("A86").Value = 43
$ws.Range("A89").Value = 44
$ws.Range("A90").Value = 45
$ws.Range("A93").Value = 46
$ws.Range("A94").Value = 47
$ws.Range("A97").Value = 48
$ws.Range("A98").Value = 49
$ws.Range("A101").Value = 50
$ws.Range("A102").Value = 51
$ws.Range("A105").Value = 52
$ws.Range("A106").Value = 53
$ws.Range("A109").Value = 54
$ws.Range("A110").Value = 55
$ws.Range("A113").Value = 56
$ws.Range("A114").Value = 57
$ws.Range("A115").Value = 58
$ws.Range("A116").Value = 59
$ws.Range("A117").Value = 60
$ws.Range("A118").Value = 61
$ws.Range("A119").Value = 62
$ws.Range("A120").Value = 63
$ws.Range("A121").Value = 64
$ws.Range("A122").Value = 65
$ws.Range("A123").Value = 66
$ws.Range("A124").Value = 67
$ws.Range("A125").Value = 68
$ws.Range("A126").Value = 69
$ws.Range("A127").Value = 70

# Methylation rate constants: 1 -> 0.001, shown with 3 decimal places
$rates = $ws.Range("E113:E127")
$rates.NumberFormat = "0.000"
$rates.Value = 0.001

# Restore the sheet selection / scroll position
$ws.Range("E113:E127").Select()
